# Corrections made during class review:
#  1. Update the cached "today" text shown in the date placeholders
#     (slide master, every slide layout and the notes master) from
#     05/05/2012 to 30/6/2012.
#  2. Shrink the title text "Atributos e métodos de uma enum" from
#     44pt to 40pt on the two slides that use it.

$p = $ppt.ActivePresentation

$oldDate = "05/05/2012"
$newDate = "30/6/2012"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        if ($shp.Type -eq 14 -or $true) {
            # Placeholder detection: guard against shapes without PlaceholderFormat
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
        }
        if ($isDatePh) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# 1a. Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster

# 1b. Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

# 1c. Notes master date placeholder.
Update-DatePlaceholder $p.NotesMaster

# 2. Shrink the "Atributos e métodos de uma enum" title on slides 6 and 7.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $txt = $shp.TextFrame.TextRange.Text
                if ($txt -eq "Atributos e métodos de uma enum") {
                    $tr = $shp.TextFrame.TextRange
                    for ($ri = 1; $ri -le $tr.Runs.Count; $ri++) {
                        $tr.Runs.Item($ri).Font.Size = 40
                    }
                    $tr.Font.Size = 40
                }
            }
        }
    }
}
